$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write plain text into a cell that does not yet contain any value,
# while forcing Excel to store it as a literal text string (avoiding the
# automatic date/number parsing that ".Value = '01/01/2018'" would trigger),
# and without generating any stray/unused cell-style entries.
function Set-PlainText($rng, $text) {
    $rng.Formula = '="' + $text + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)  # xlPasteValues
    $rng.Worksheet.Application.CutCopyMode = $false
}

function Copy-CellFormat($src, $dst) {
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
    $src.Worksheet.Application.CutCopyMode = $false
}

# The old row 13 only held a stray "Docentes responsáveis" value in B/C with
# no label in column A; it is removed entirely and everything below shifts
# up by one row, carrying its row height/style formatting along with it.
$ws.Rows(13).Delete()

# Row 10 (Objetivos:) now shows the "Docentes responsáveis" value instead of
# the long objectives paragraph.
$ws.Range("B10").Value = "5840712 - Ângelo Capri Neto"
$ws.Range("C10").Value = "5840712 - Ângelo Capri Neto"

# Row 13 (Programa resumido:) now shows "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) now shows the activation date as plain text. B15/C15
# were empty before, so copy the text-cell formatting from B8/C8 first.
Copy-CellFormat $ws.Range("B8") $ws.Range("B15")
Set-PlainText $ws.Range("B15") "01/01/2018"

Copy-CellFormat $ws.Range("C8") $ws.Range("C15")
Set-PlainText $ws.Range("C15") "01/01/2018"

# Row 18 (Método:) now shows the "Docentes responsáveis" value.
$ws.Range("B18").Value = "5840712 - Ângelo Capri Neto"
$ws.Range("C18").Value = "5840712 - Ângelo Capri Neto"

# Row 19 (Critério:) now shows the evaluation method text.
$metodoText = "A avaliação será feita por meio de duas provas (P1 e P2). A critério do professor, a avaliação poderá ser complementada por meio de trabalhos e/ou relatórios, valendo até 30% da nota das provas."
$ws.Range("B19").Value = $metodoText
$ws.Range("C19").Value = $metodoText

# Row 20 (Norma de recuperação:) now shows the final grade formula text.
$notaFinalText = "A nota final (NF) será calculada pela média aritmética das provas. NF=(P1 +P2)/2."
$ws.Range("B20").Value = $notaFinalText
$ws.Range("C20").Value = $notaFinalText

# Row 21 (Bibliografia:) now shows the recovery exam rule text.
$recuperacaoText = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2."
$ws.Range("B21").Value = $recuperacaoText
$ws.Range("C21").Value = $recuperacaoText
